# Controls.xlsx update — "Added a new control" (Konteringsmallar, Momskontroll,
# Periodiseringar, Fakturahantering sections) per the module-fetch log cleanup
# described in the commit message. New rows are appended after the existing
# data (through row 455) and written in the same order the originating script
# emitted its log lines, which is why the resulting shared-string table is
# populated sequentially even though the target row numbers are not strictly
# increasing in that same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ControlRow {
    param(
        [int]$Row,
        [string]$WindowClassName,
        [int]$ControlId,
        [string]$Module,
        [string]$Text
    )

    $ws.Cells.Item($Row, 1).Value2 = $WindowClassName
    $ws.Cells.Item($Row, 2).Value2 = $ControlId
    $ws.Cells.Item($Row, 3).Value2 = $Module
    $ws.Cells.Item($Row, 4).Value2 = $Text
}

# Konteringsmallar (accounting templates) grid + calculator-style buttons.
Set-ControlRow 456 "SafGrid" 20723 "Konteringsmallar" "Konteringsmallar"
Set-ControlRow 458 "Button"  21146 "Konteringsmallar" "Minus"
Set-ControlRow 457 "Button"  21147 "Konteringsmallar" "Plus"
Set-ControlRow 460 "Button"  21149 "Konteringsmallar" "Dividera"
Set-ControlRow 459 "Button"  21148 "Konteringsmallar" "Multiplicera"
Set-ControlRow 461 "Button"  21150 "Konteringsmallar" "Vänster parentes"
Set-ControlRow 462 "Button"  21151 "Konteringsmallar" "Höger parentes"
Set-ControlRow 463 "Button"  23637 "Konteringsmallar" "Vänster klammer"
Set-ControlRow 464 "Button"  23638 "Konteringsmallar" "Höger klammer"

# Gamla journaler / Journaler.
Set-ControlRow 465 "SafGrid" 20723 "Gamla journaler" "Journaler"

# Momskontroll (VAT check) date range edits.
Set-ControlRow 467 "Edit" 22085 "Momskontroll" "Till och med"
Set-ControlRow 466 "Edit" 22084 "Momskontroll" "Från och med"

# Periodiseringar (accruals) list.
Set-ControlRow 468 "SafGrid" 24908 "Periodiseringar" "Lista över periodiseringar"

# Fakturahantering (invoice handling) list.
Set-ControlRow 469 "SafGrid" 21346 "Fakturahantering" "Fakturor"

# Leave the selection where the author ended up after entering the last row.
$ws.Range("D469").Select() | Out-Null
